$d = $word.ActiveDocument

# --- Merge the split date-range runs into a single run: "Mar 2022 - current" ---
# The runs "Mar" / " 202" / "2" / " " / "-" / " " / "current" already concatenate
# to the exact desired text, so a Find/Replace of the full phrase with itself
# collapses them into one run without altering the visible content.
$dateRange = $d.Content
$dateRange.Find.Execute("Mar 2022 - current", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Mar 2022 - current", 2)

# --- Fix typo: "Healthy Analysis Division" => "Health Analysis Division" ---
# Both occurrences live inside a hyperlink run (rStyle="Hyperlink"). Rather than
# using Find/Replace (which rebuilds the run and drops rStyle/formatting here),
# surgically delete just the stray "y" in "Healthy" so the run keeps its
# original run properties and is merely shortened by one character.
$continueSearch = $true
while ($continueSearch) {
    $findRange = $d.Content
    $findRange.Find.Execute("Healthy Analysis Division")
    if ($findRange.Find.Found) {
        $yStart = $findRange.Start + 6
        $yRange = $d.Range($yStart, $yStart + 1)
        $yRange.Delete()
    } else {
        $continueSearch = $false
    }
}
